# Fruta / hortaliza, semanal
#
# Two new weekly price records (2023-04-28) are inserted into the
# consolidated "Pepino dulce" sheet right after the existing 2021-07-02
# entries (row 30), pushing every subsequent row down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new records.
$ws.Rows(31).Insert()
$ws.Rows(31).Insert()

# Row 31 - "Primera" quality
$ws.Cells.Item(31, 1).Value  = 11
$ws.Cells.Item(31, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(31, 3).Value  = "Bíobío"
$ws.Cells.Item(31, 4).Value  = [DateTime]"2023-04-28"
$ws.Cells.Item(31, 5).Value  = 8
$ws.Cells.Item(31, 6).Value  = 100112043
$ws.Cells.Item(31, 7).Value  = "Pepino dulce"
$ws.Cells.Item(31, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(31, 9).Value  = "Primera"
$ws.Cells.Item(31, 10).Value = 220
$ws.Cells.Item(31, 11).Value = 11000
$ws.Cells.Item(31, 12).Value = 12000
$ws.Cells.Item(31, 13).Value = 11545
$ws.Cells.Item(31, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(31, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value = 641
$ws.Cells.Item(31, 17).Value = 18
$ws.Cells.Item(31, 18).Value = "Hortaliza"

# Row 32 - "Segunda" quality
$ws.Cells.Item(32, 1).Value  = 11
$ws.Cells.Item(32, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(32, 3).Value  = "Bíobío"
$ws.Cells.Item(32, 4).Value  = [DateTime]"2023-04-28"
$ws.Cells.Item(32, 5).Value  = 8
$ws.Cells.Item(32, 6).Value  = 100112043
$ws.Cells.Item(32, 7).Value  = "Pepino dulce"
$ws.Cells.Item(32, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(32, 9).Value  = "Segunda"
$ws.Cells.Item(32, 10).Value = 100
$ws.Cells.Item(32, 11).Value = 10000
$ws.Cells.Item(32, 12).Value = 10000
$ws.Cells.Item(32, 13).Value = 10000
$ws.Cells.Item(32, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(32, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(32, 16).Value = 556
$ws.Cells.Item(32, 17).Value = 18
$ws.Cells.Item(32, 18).Value = "Hortaliza"
